# New crime data collected - update the CompStat weekly report figures.
#
# 1) Update the report header text (volume/issue number and the
#    week-covering date range) stored as shared strings.
# 2) Update the Crime Complaints table (rows 14-31, "TOTAL" row 21,
#    and the historical "14 Year" row 33) with the newly collected
#    weekly figures and their recalculated percentage changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
# "Volume 31   Number  44" -> "Volume 31   Number  45"
$ws.Range("A8").Value = "Volume 31   Number  45"
# "Report Covering the Week  10/28/2024  Through  11/3/2024"
#   -> "Report Covering the Week  11/4/2024  Through  11/10/2024"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Crime Complaints table ------------------------------------------
$ws.Range("F14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("I14").Value = 18
$ws.Range("K14").Value = 20
$ws.Range("L14").Value = -33.333333333333
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = -75.675675675675
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 150
$ws.Range("F15").Value = 16
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 141
$ws.Range("J15").Value = 101
$ws.Range("K15").Value = 39.603960396039
$ws.Range("L15").Value = -10.191082802547
$ws.Range("M15").Value = 46.875
$ws.Range("N15").Value = -25.789473684210
$ws.Range("C16").Value = 44
$ws.Range("D16").Value = 42
$ws.Range("E16").Value = 4.761904761904
$ws.Range("F16").Value = 150
$ws.Range("G16").Value = 134
$ws.Range("H16").Value = 11.940298507462
$ws.Range("I16").Value = 1484
$ws.Range("J16").Value = 1550
$ws.Range("K16").Value = -4.258064516129
$ws.Range("L16").Value = -20.300751879699
$ws.Range("M16").Value = 29.833770778652
$ws.Range("N16").Value = -83.733421023786
$ws.Range("C17").Value = 41
$ws.Range("D17").Value = 40
$ws.Range("E17").Value = 2.5
$ws.Range("F17").Value = 183
$ws.Range("G17").Value = 138
$ws.Range("H17").Value = 32.608695652173
$ws.Range("I17").Value = 2049
$ws.Range("J17").Value = 1874
$ws.Range("K17").Value = 9.338313767342
$ws.Range("L17").Value = 10.637149028077
$ws.Range("M17").Value = 91.674462114125
$ws.Range("N17").Value = -32.152317880794
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 63
$ws.Range("E18").Value = -53.968253968254
$ws.Range("F18").Value = 118
$ws.Range("G18").Value = 199
$ws.Range("H18").Value = -40.703517587939
$ws.Range("I18").Value = 1560
$ws.Range("J18").Value = 1873
$ws.Range("K18").Value = -16.711158569140
$ws.Range("L18").Value = -39.417475728155
$ws.Range("M18").Value = 1.430429128738
$ws.Range("N18").Value = -85.058902403984
$ws.Range("C19").Value = 207
$ws.Range("D19").Value = 244
$ws.Range("E19").Value = -15.163934426229
$ws.Range("F19").Value = 889
$ws.Range("G19").Value = 990
$ws.Range("H19").Value = -10.202020202020
$ws.Range("I19").Value = 9092
$ws.Range("J19").Value = 10042
$ws.Range("K19").Value = -9.460266879107
$ws.Range("L19").Value = -10.503002264002
$ws.Range("M19").Value = -0.153744783659
$ws.Range("N19").Value = -68.145189545231
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 11.111111111111
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = -3.225806451612
$ws.Range("I20").Value = 386
$ws.Range("J20").Value = 542
$ws.Range("K20").Value = -28.782287822878
$ws.Range("L20").Value = -35.451505016722
$ws.Range("M20").Value = 10.601719197707
$ws.Range("N20").Value = -92.887414777962
$ws.Range("C21").Value = 337
$ws.Range("D21").Value = 400
$ws.Range("E21").Value = -15.75
$ws.Range("F21").Value = 1387
$ws.Range("G21").Value = 1501
$ws.Range("H21").Value = -7.594936708860
$ws.Range("I21").Value = 14730
$ws.Range("J21").Value = 15997
$ws.Range("K21").Value = -7.920235044070
$ws.Range("L21").Value = -14.509576320371
$ws.Range("M21").Value = 10.618804445779
$ws.Range("N21").Value = -74.074660752943
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = 17
$ws.Range("E22").Value = -11.764705882352
$ws.Range("F22").Value = 52
$ws.Range("G22").Value = 62
$ws.Range("H22").Value = -16.129032258064
$ws.Range("I22").Value = 524
$ws.Range("J22").Value = 585
$ws.Range("K22").Value = -10.427350427350
$ws.Range("L22").Value = -10.273972602739
$ws.Range("M22").Value = 12.688172043010
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = -63.636363636363
$ws.Range("F23").Value = 36
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = 12.5
$ws.Range("I23").Value = 350
$ws.Range("J23").Value = 348
$ws.Range("K23").Value = 0.574712643678
$ws.Range("L23").Value = -8.616187989556
$ws.Range("M23").Value = 13.268608414239
$ws.Range("C24").Value = 485
$ws.Range("D24").Value = 458
$ws.Range("E24").Value = 5.895196506550
$ws.Range("F24").Value = 1796
$ws.Range("G24").Value = 1787
$ws.Range("H24").Value = 0.503637381085
$ws.Range("I24").Value = 19344
$ws.Range("J24").Value = 18072
$ws.Range("K24").Value = 7.038512616201
$ws.Range("L24").Value = 2.376290023815
$ws.Range("M24").Value = 32.366224168605
$ws.Range("C25").Value = 403
$ws.Range("D25").Value = 380
$ws.Range("E25").Value = 6.052631578947
$ws.Range("F25").Value = 1442
$ws.Range("G25").Value = 1471
$ws.Range("H25").Value = -1.971447994561
$ws.Range("I25").Value = 16215
$ws.Range("J25").Value = 14916
$ws.Range("K25").Value = 8.708769106999
$ws.Range("L25").Value = 1.135158735108
$ws.Range("C26").Value = 89
$ws.Range("D26").Value = 93
$ws.Range("E26").Value = -4.301075268817
$ws.Range("F26").Value = 454
$ws.Range("G26").Value = 430
$ws.Range("H26").Value = 5.581395348837
$ws.Range("I26").Value = 4476
$ws.Range("J26").Value = 4355
$ws.Range("K26").Value = 2.778415614236
$ws.Range("L26").Value = 9.571603427172
$ws.Range("M26").Value = 42.638623326959
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 25
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 206
$ws.Range("J27").Value = 177
$ws.Range("K27").Value = 16.384180790960
$ws.Range("L27").Value = -15.918367346938
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 17
$ws.Range("E28").Value = 29.411764705882
$ws.Range("F28").Value = 73
$ws.Range("G28").Value = 85
$ws.Range("H28").Value = -14.117647058823
$ws.Range("I28").Value = 835
$ws.Range("J28").Value = 791
$ws.Range("K28").Value = 5.562579013906
$ws.Range("L28").Value = -2.110199296600
$ws.Range("D29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("J29").Value = 33
$ws.Range("K29").Value = 6.060606060606
$ws.Range("N29").Value = -70.338983050847
$ws.Range("D30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 28
$ws.Range("K30").Value = 7.142857142857
$ws.Range("N30").Value = -71.698113207547
$ws.Range("D31").Value = 3
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 24
$ws.Range("H31").Value = -83.333333333333
$ws.Range("J31").Value = 121
$ws.Range("K31").Value = 3.305785123966
$ws.Range("L31").Value = -13.194444444444
$ws.Range("D33").Copy($ws.Range("C33"))
